$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the monthly columns right by one to make room for the new month,
# then set the new month in column D. This turns:
#   D1=Dec_2025, E1=Nov_2025, F1=Oct_2025
# into:
#   D1=Jan_2026, E1=Dec_2025, F1=Nov_2025
# (Oct_2025 falls off and is discarded; G1/MoM and H1/QoQ are untouched)

$ws.Range("F1").Value = $ws.Range("E1").Value2
$ws.Range("E1").Value = $ws.Range("D1").Value2
$ws.Range("D1").Value = "Jan_2026"
